# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newly generated data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1729
    $ws.Range("F3").Value = 7947
    $ws.Range("F4").Value = 184
    $ws.Range("F5").Value = 267
}
